$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=STUDYID, B=DOMAIN, C=IETESTCD, D=IETEST, E=IECAT, F=IESCAT, G=IEORRES
# Row 1 = header (unchanged). Data rows: inclusion criteria in rows 2-14,
# exclusion criteria in rows 15-30 (replacing/expanding the original rows 2-19).

$ws.Cells.Item(2, 1).Value = 'STUDY001'
$ws.Cells.Item(2, 2).Value = 'TI'
$ws.Cells.Item(2, 3).Value = 'INCL001'
$ws.Cells.Item(2, 4).Value = 'Inclusion Criteria'
$ws.Cells.Item(2, 5).Value = ''
$ws.Cells.Item(2, 6).Value = ''
$ws.Cells.Item(2, 7).Value = 'HER2-positive breast cancer HER2-positive status will be based on pretreatment biopsy material and defined as an immunohistochemistry (IHC) (Appendix 6) score of 3+ and/or posit... (As per the protocol)'

$ws.Cells.Item(3, 1).Value = 'STUDY001'
$ws.Cells.Item(3, 2).Value = 'TI'
$ws.Cells.Item(3, 3).Value = 'INCL002'
$ws.Cells.Item(3, 4).Value = 'Inclusion Criteria'
$ws.Cells.Item(3, 5).Value = ''
$ws.Cells.Item(3, 6).Value = ''
$ws.Cells.Item(3, 7).Value = 'Histologically confirmed invasive breast carcinoma'

$ws.Cells.Item(4, 1).Value = 'STUDY001'
$ws.Cells.Item(4, 2).Value = 'TI'
$ws.Cells.Item(4, 3).Value = 'INCL003'
$ws.Cells.Item(4, 4).Value = 'Inclusion Criteria'
$ws.Cells.Item(4, 5).Value = ''
$ws.Cells.Item(4, 6).Value = ''
$ws.Cells.Item(4, 7).Value = 'Clinical stage at presentation: T1–4, N0–3, M0 (Note: Patients with T1a/bN0 tumors will not be eligible)'

$ws.Cells.Item(5, 1).Value = 'STUDY001'
$ws.Cells.Item(5, 2).Value = 'TI'
$ws.Cells.Item(5, 3).Value = 'INCL004'
$ws.Cells.Item(5, 4).Value = 'Inclusion Criteria'
$ws.Cells.Item(5, 5).Value = ''
$ws.Cells.Item(5, 6).Value = ''
$ws.Cells.Item(5, 7).Value = 'Completion of preoperative systemic chemotherapy and HER2-directed treatment. Systemic therapy must consist of at least 6 cycles of chemotherapy, with a total duration at least... (As per the protocol)'

$ws.Cells.Item(6, 1).Value = 'STUDY001'
$ws.Cells.Item(6, 2).Value = 'TI'
$ws.Cells.Item(6, 3).Value = 'INCL005'
$ws.Cells.Item(6, 4).Value = 'Inclusion Criteria'
$ws.Cells.Item(6, 5).Value = ''
$ws.Cells.Item(6, 6).Value = ''
$ws.Cells.Item(6, 7).Value = 'Adequate excision: surgical removal of all clinically evident disease in the breast and lymph nodes as follows: Breast surgery: total mastectomy with no gross residual disease a... (As per the protocol)'

$ws.Cells.Item(7, 1).Value = 'STUDY001'
$ws.Cells.Item(7, 2).Value = 'TI'
$ws.Cells.Item(7, 3).Value = 'INCL006'
$ws.Cells.Item(7, 4).Value = 'Inclusion Criteria'
$ws.Cells.Item(7, 5).Value = ''
$ws.Cells.Item(7, 6).Value = ''
$ws.Cells.Item(7, 7).Value = 'Pathologic evidence of residual invasive carcinoma in the breast or axillary lymph nodes following completion of preoperative therapy. If invasive disease is present in both bre... (As per the protocol)'

$ws.Cells.Item(8, 1).Value = 'STUDY001'
$ws.Cells.Item(8, 2).Value = 'TI'
$ws.Cells.Item(8, 3).Value = 'INCL007'
$ws.Cells.Item(8, 4).Value = 'Inclusion Criteria'
$ws.Cells.Item(8, 5).Value = ''
$ws.Cells.Item(8, 6).Value = ''
$ws.Cells.Item(8, 7).Value = 'An interval of no more than 12 weeks between the date of primary surgery and the date of randomization'

$ws.Cells.Item(9, 1).Value = 'STUDY001'
$ws.Cells.Item(9, 2).Value = 'TI'
$ws.Cells.Item(9, 3).Value = 'INCL008'
$ws.Cells.Item(9, 4).Value = 'Inclusion Criteria'
$ws.Cells.Item(9, 5).Value = ''
$ws.Cells.Item(9, 6).Value = ''
$ws.Cells.Item(9, 7).Value = 'Known hormone receptor status Hormone receptor−positive status can be determined by either known positive ER or known positive PgR status; hormone receptor−negative status must... (As per the protocol)'

$ws.Cells.Item(10, 1).Value = 'STUDY001'
$ws.Cells.Item(10, 2).Value = 'TI'
$ws.Cells.Item(10, 3).Value = 'INCL009'
$ws.Cells.Item(10, 4).Value = 'Inclusion Criteria'
$ws.Cells.Item(10, 5).Value = ''
$ws.Cells.Item(10, 6).Value = ''
$ws.Cells.Item(10, 7).Value = 'Signed written informed consent approved by the study site’s Institutional Review Board (IRB)/Ethical Committee (EC)'

$ws.Cells.Item(11, 1).Value = 'STUDY001'
$ws.Cells.Item(11, 2).Value = 'TI'
$ws.Cells.Item(11, 3).Value = 'INCL010'
$ws.Cells.Item(11, 4).Value = 'Inclusion Criteria'
$ws.Cells.Item(11, 5).Value = ''
$ws.Cells.Item(11, 6).Value = ''
$ws.Cells.Item(11, 7).Value = 'Age ≥ 18 years'

$ws.Cells.Item(12, 1).Value = 'STUDY001'
$ws.Cells.Item(12, 2).Value = 'TI'
$ws.Cells.Item(12, 3).Value = 'INCL011'
$ws.Cells.Item(12, 4).Value = 'Inclusion Criteria'
$ws.Cells.Item(12, 5).Value = ''
$ws.Cells.Item(12, 6).Value = ''
$ws.Cells.Item(12, 7).Value = 'Eastern Cooperative Oncology Group (ECOG) performance status 0 or 1'

$ws.Cells.Item(13, 1).Value = 'STUDY001'
$ws.Cells.Item(13, 2).Value = 'TI'
$ws.Cells.Item(13, 3).Value = 'INCL012'
$ws.Cells.Item(13, 4).Value = 'Inclusion Criteria'
$ws.Cells.Item(13, 5).Value = ''
$ws.Cells.Item(13, 6).Value = ''
$ws.Cells.Item(13, 7).Value = 'Life expectancy ≥ 6 months'

$ws.Cells.Item(14, 1).Value = 'STUDY001'
$ws.Cells.Item(14, 2).Value = 'TI'
$ws.Cells.Item(14, 3).Value = 'INCL013'
$ws.Cells.Item(14, 4).Value = 'Inclusion Criteria'
$ws.Cells.Item(14, 5).Value = ''
$ws.Cells.Item(14, 6).Value = ''
$ws.Cells.Item(14, 7).Value = 'Adequate organ function during screening, defined as: a. Absolute neutrophil count ≥ 1200 cells/mm3 b. Platelet count ≥ 100000 cells/mm3 c. Hemoglobin ≥ 9.0 g/dL; patients may r... (As per the protocol)'

$ws.Cells.Item(15, 1).Value = 'STUDY001'
$ws.Cells.Item(15, 2).Value = 'TI'
$ws.Cells.Item(15, 3).Value = 'EXCL001'
$ws.Cells.Item(15, 4).Value = 'Exclusion Criteria'
$ws.Cells.Item(15, 5).Value = ''
$ws.Cells.Item(15, 6).Value = ''
$ws.Cells.Item(15, 7).Value = 'Stage IV (metastatic) breast cancer'

$ws.Cells.Item(16, 1).Value = 'STUDY001'
$ws.Cells.Item(16, 2).Value = 'TI'
$ws.Cells.Item(16, 3).Value = 'EXCL002'
$ws.Cells.Item(16, 4).Value = 'Exclusion Criteria'
$ws.Cells.Item(16, 5).Value = ''
$ws.Cells.Item(16, 6).Value = ''
$ws.Cells.Item(16, 7).Value = 'History of any prior (ipsi- or contralateral) breast cancer except lobular CIS'

$ws.Cells.Item(17, 1).Value = 'STUDY001'
$ws.Cells.Item(17, 2).Value = 'TI'
$ws.Cells.Item(17, 3).Value = 'EXCL003'
$ws.Cells.Item(17, 4).Value = 'Exclusion Criteria'
$ws.Cells.Item(17, 5).Value = ''
$ws.Cells.Item(17, 6).Value = ''
$ws.Cells.Item(17, 7).Value = 'Evidence of clinically evident gross residual or recurrent disease following preoperative therapy and surgery'

$ws.Cells.Item(18, 1).Value = 'STUDY001'
$ws.Cells.Item(18, 2).Value = 'TI'
$ws.Cells.Item(18, 3).Value = 'EXCL004'
$ws.Cells.Item(18, 4).Value = 'Exclusion Criteria'
$ws.Cells.Item(18, 5).Value = ''
$ws.Cells.Item(18, 6).Value = ''
$ws.Cells.Item(18, 7).Value = 'An overall response of PD according to the investigator at the conclusion of preoperative systemic therapy'

$ws.Cells.Item(19, 1).Value = 'STUDY001'
$ws.Cells.Item(19, 2).Value = 'TI'
$ws.Cells.Item(19, 3).Value = 'EXCL005'
$ws.Cells.Item(19, 4).Value = 'Exclusion Criteria'
$ws.Cells.Item(19, 5).Value = ''
$ws.Cells.Item(19, 6).Value = ''
$ws.Cells.Item(19, 7).Value = 'Treatment with any anti-cancer investigational drug within 28 days prior to commencing study treatment Trastuzumab Emtansine—F. Hoffmann-La Roche Ltd 43'

$ws.Cells.Item(20, 1).Value = 'STUDY001'
$ws.Cells.Item(20, 2).Value = 'TI'
$ws.Cells.Item(20, 3).Value = 'EXCL006'
$ws.Cells.Item(20, 4).Value = 'Exclusion Criteria'
$ws.Cells.Item(20, 5).Value = ''
$ws.Cells.Item(20, 6).Value = ''
$ws.Cells.Item(20, 7).Value = 'History of other malignancy within the last 5 years except for appropriately treated CIS of the cervix, non-melanoma skin carcinoma, Stage I uterine cancer, or other non-breast... (As per the protocol)'

$ws.Cells.Item(21, 1).Value = 'STUDY001'
$ws.Cells.Item(21, 2).Value = 'TI'
$ws.Cells.Item(21, 3).Value = 'EXCL007'
$ws.Cells.Item(21, 4).Value = 'Exclusion Criteria'
$ws.Cells.Item(21, 5).Value = ''
$ws.Cells.Item(21, 6).Value = ''
$ws.Cells.Item(21, 7).Value = 'Patients for whom radiotherapy would be recommended for breast cancer treatment but for whom it is contraindicated because of medical reasons (e.g., connective tissue disorder o... (As per the protocol)'

$ws.Cells.Item(22, 1).Value = 'STUDY001'
$ws.Cells.Item(22, 2).Value = 'TI'
$ws.Cells.Item(22, 3).Value = 'EXCL008'
$ws.Cells.Item(22, 4).Value = 'Exclusion Criteria'
$ws.Cells.Item(22, 5).Value = ''
$ws.Cells.Item(22, 6).Value = ''
$ws.Cells.Item(22, 7).Value = 'Current NCI CTCAE (Version 4.0) Grade ≥ 2 peripheral neuropathy'

$ws.Cells.Item(23, 1).Value = 'STUDY001'
$ws.Cells.Item(23, 2).Value = 'TI'
$ws.Cells.Item(23, 3).Value = 'EXCL009'
$ws.Cells.Item(23, 4).Value = 'Exclusion Criteria'
$ws.Cells.Item(23, 5).Value = ''
$ws.Cells.Item(23, 6).Value = ''
$ws.Cells.Item(23, 7).Value = 'History of exposure to the following cumulative doses of anthracyclines: Doxorubicin > 240 mg/m2 Epirubicin or Liposomal Doxorubicin-Hydrochloride (Myocet®) > 480 mg/m2 For othe... (As per the protocol)'

$ws.Cells.Item(24, 1).Value = 'STUDY001'
$ws.Cells.Item(24, 2).Value = 'TI'
$ws.Cells.Item(24, 3).Value = 'EXCL010'
$ws.Cells.Item(24, 4).Value = 'Exclusion Criteria'
$ws.Cells.Item(24, 5).Value = ''
$ws.Cells.Item(24, 6).Value = ''
$ws.Cells.Item(24, 7).Value = 'Cardiopulmonary dysfunction as defined by any of the following: History of NCI CTCAE (Version 4.0) Grade ≥ 3 symptomatic CHF or NYHA criteria Class ≥ II Angina pectoris requirin... (As per the protocol)'

$ws.Cells.Item(25, 1).Value = 'STUDY001'
$ws.Cells.Item(25, 2).Value = 'TI'
$ws.Cells.Item(25, 3).Value = 'EXCL011'
$ws.Cells.Item(25, 4).Value = 'Exclusion Criteria'
$ws.Cells.Item(25, 5).Value = ''
$ws.Cells.Item(25, 6).Value = ''
$ws.Cells.Item(25, 7).Value = 'High-risk uncontrolled arrhythmias: i.e., atrial tachycardia with a heart rate > 100/min at rest, significant ventricular arrhythmia (ventricular tachycardia) or higher-grade AV... (As per the protocol)'

$ws.Cells.Item(26, 1).Value = 'STUDY001'
$ws.Cells.Item(26, 2).Value = 'TI'
$ws.Cells.Item(26, 3).Value = 'EXCL012'
$ws.Cells.Item(26, 4).Value = 'Exclusion Criteria'
$ws.Cells.Item(26, 5).Value = ''
$ws.Cells.Item(26, 6).Value = ''
$ws.Cells.Item(26, 7).Value = 'Prior treatment with trastuzumab emtansine'

$ws.Cells.Item(27, 1).Value = 'STUDY001'
$ws.Cells.Item(27, 2).Value = 'TI'
$ws.Cells.Item(27, 3).Value = 'EXCL013'
$ws.Cells.Item(27, 4).Value = 'Exclusion Criteria'
$ws.Cells.Item(27, 5).Value = ''
$ws.Cells.Item(27, 6).Value = ''
$ws.Cells.Item(27, 7).Value = 'Current severe, uncontrolled systemic disease (e.g., clinically significant cardiovascular, pulmonary, or metabolic disease; wound-healing disorders; ulcers)'

$ws.Cells.Item(28, 1).Value = 'STUDY001'
$ws.Cells.Item(28, 2).Value = 'TI'
$ws.Cells.Item(28, 3).Value = 'EXCL014'
$ws.Cells.Item(28, 4).Value = 'Exclusion Criteria'
$ws.Cells.Item(28, 5).Value = ''
$ws.Cells.Item(28, 6).Value = ''
$ws.Cells.Item(28, 7).Value = 'For female patients, current pregnancy and/or lactation'

$ws.Cells.Item(29, 1).Value = 'STUDY001'
$ws.Cells.Item(29, 2).Value = 'TI'
$ws.Cells.Item(29, 3).Value = 'EXCL015'
$ws.Cells.Item(29, 4).Value = 'Exclusion Criteria'
$ws.Cells.Item(29, 5).Value = ''
$ws.Cells.Item(29, 6).Value = ''
$ws.Cells.Item(29, 7).Value = 'Major surgical procedure unrelated to breast cancer or significant traumatic injury within approximately 28 days prior to randomization or anticipation of the need for major sur... (As per the protocol)'

$ws.Cells.Item(30, 1).Value = 'STUDY001'
$ws.Cells.Item(30, 2).Value = 'TI'
$ws.Cells.Item(30, 3).Value = 'EXCL016'
$ws.Cells.Item(30, 4).Value = 'Exclusion Criteria'
$ws.Cells.Item(30, 5).Value = ''
$ws.Cells.Item(30, 6).Value = ''
$ws.Cells.Item(30, 7).Value = 'Any known active liver disease, for example, disease due to HBV, HCV, autoimmune hepatic disorders, or sclerosing cholangitis. Patients who have positive Trastuzumab Emtansine—F... (As per the protocol)'
